# Trade #13 closed at 2026-02-17 13:10:06 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = -0.62
$wsSummary.Range("B6").Value = 13
$wsSummary.Range("B9").Value = 30.77

# --- Strategy Status sheet ---
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("D4").Value = 13
$wsStrategy.Range("G4").Value = 30.77

# --- New trade row data (Trade #13) ---
$tradeRow = @(13, "2026-02-17", "13:08:55", "MarketMaking", "DOWN", 0.98, 0.98, "CLOSED", 0, 0, 99.59999999999999, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 1.23)

# --- All Trades sheet: append new row 14 ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
# Force column B (date-looking text) to be stored as text, not auto-parsed into a date serial
$wsAllTrades.Range("B14").NumberFormat = "@"
for ($i = 0; $i -lt $tradeRow.Length; $i++) {
    $wsAllTrades.Cells.Item(14, $i + 1).Value = $tradeRow[$i]
}

# --- MarketMaking sheet: append new row 14 ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
$wsMarketMaking.Range("B14").NumberFormat = "@"
for ($i = 0; $i -lt $tradeRow.Length; $i++) {
    $wsMarketMaking.Cells.Item(14, $i + 1).Value = $tradeRow[$i]
}
